$wb = $excel.ActiveWorkbook

# --- Sheet "Voucher": change B10 from 495.0 to 0.0 ---
$voucher = $wb.Worksheets.Item("Voucher")
$voucher.Range("B10").Value = 0.0

# --- Sheet "Upload": delete row 6 (the "Good Faith Acceptance Credit" row),
#     which shifts all subsequent rows up by one and removes the old last row (13) ---
$upload = $wb.Worksheets.Item("Upload")
$upload.Rows.Item(6).Delete()
